$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $Text)
    $cell = $ws.Range($CellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "29.186.22"
$ws.Range("E2").Value = "  -2.81%  "
$ws.Range("D3").Value = "1.850.06"
$ws.Range("E3").Value = "  -1.88%  "
Set-TextValue "D4" "1.000"
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue "D5" "0.7043"
$ws.Range("E5").Value = "  -4.31%  "
Set-TextValue "D6" "238.93"
$ws.Range("E6").Value = "  -1.27%  "
$ws.Range("E7").Value = "  +0.01%  "
Set-TextValue "D8" "0.3055"
$ws.Range("E8").Value = "  -3.47%  "
Set-TextValue "D9" "0.07446"
$ws.Range("E9").Value = "  +3.63%  "
Set-TextValue "D10" "23.46"
$ws.Range("E10").Value = "  -5.19%  "
Set-TextValue "D11" "0.08149"
$ws.Range("E11").Value = "  -2.10%  "
$ws.Range("D12").Value = "1.907.02"
$ws.Range("E12").Value = "  +0.70%  "
Set-TextValue "D13" "0.7291"
$ws.Range("E13").Value = "  -3.66%  "
Set-TextValue "D14" "5.221"
$ws.Range("E14").Value = "  -3.26%  "
$ws.Range("E15").Value = "  -4.00%  "
$ws.Range("D16").Value = "29.274.37"
$ws.Range("E16").Value = "  -2.55%  "
Set-TextValue "D17" "5.792"
$ws.Range("E17").Value = "  -6.03%  "
Set-TextValue "D18" "239.03"
$ws.Range("E18").Value = "  -4.38%  "
Set-TextValue "D19" "13.12"
Set-TextValue "D20" "0.000007659"
$ws.Range("E20").Value = "  -2.49%  "
Set-TextValue "D21" "1.000"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").Value = "2.126.02"
$ws.Range("E22").Value = "  -1.63%  "
$ws.Range("E23").Value = "  +0.10%  "
Set-TextValue "D24" "7.603"
$ws.Range("E24").Value = "  -3.70%  "
Set-TextValue "D25" "9.026"
$ws.Range("E25").Value = "  -2.68%  "
Set-TextValue "D26" "0.1464"
$ws.Range("E26").Value = "  -6.19%  "
Set-TextValue "D27" "161.19"
$ws.Range("E27").Value = "  -1.21%  "
Set-TextValue "D28" "18.13"
$ws.Range("E28").Value = "  -2.92%  "
Set-TextValue "D29" "1.975"
$ws.Range("E29").Value = "  -3.59%  "
Set-TextValue "D30" "1.411"
$ws.Range("E30").Value = "  -4.50%  "
Set-TextValue "D31" "4.512"
$ws.Range("E31").Value = "  -1.19%  "
$ws.Range("E32").Value = "  -2.64%  "
Set-TextValue "D33" "4.013"
$ws.Range("E33").Value = "  -4.47%  "
Set-TextValue "D34" "0.05207"
$ws.Range("E34").Value = "  -2.40%  "
Set-TextValue "D35" "1.190"
$ws.Range("E35").Value = "  -4.65%  "
Set-TextValue "D36" "1.045"
$ws.Range("E36").Value = "  +4.68%  "
Set-TextValue "D37" "0.7072"
$ws.Range("E37").Value = "  -8.04%  "
Set-TextValue "D38" "2.657"
$ws.Range("E38").Value = "  -2.26%  "
Set-TextValue "D39" "0.01871"
$ws.Range("E39").Value = "  -4.53%  "
Set-TextValue "D40" "2.681"
$ws.Range("E40").Value = "  -2.84%  "
Set-TextValue "D41" "0.9426"
$ws.Range("E41").Value = "  +8.02%  "
Set-TextValue "D42" "6.031"
$ws.Range("E42").Value = "  -0.15%  "
Set-TextValue "D43" "0.4310"
$ws.Range("E43").Value = "  -5.80%  "
$ws.Range("D44").Value = "1.067.37"
$ws.Range("E44").Value = "  -1.93%  "
Set-TextValue "D45" "70.64"
$ws.Range("E45").Value = "  -2.29%  "
Set-TextValue "D46" "0.9999"
$ws.Range("E46").Value = "  +0.03%  "
Set-TextValue "D47" "103.61"
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("D48").Value = "2.027.75"
$ws.Range("E48").Value = "  -1.08%  "
Set-TextValue "D49" "1.750"
$ws.Range("E49").Value = "  -5.76%  "
Set-TextValue "D50" "7.058"
$ws.Range("E50").Value = "  -6.84%  "
Set-TextValue "D51" "9.105"
$ws.Range("E51").Value = "  -4.74%  "
